$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The base "servings" quantity for this recipe dropped from 4 to 3 for every
# ingredient row (E2:E8). The dependent "Daily recipe" column (B) is a
# formula (=E2*Fx/Gx, shared for B3:B8) and recalculates automatically.
$ws.Range("E2:E8").Value = 3

# The gram quantity for the last ingredient (row 8) also changed.
$ws.Range("F8").Value = 2

# Remove the now-unused "water in first go" / "water in second go" helper
# rows entirely (label + formula).
$ws.Range("A10:B11").ClearContents()

# Update the selected cell shown when the workbook was last saved.
$ws.Range("F8").Select()

$excel.Calculate()
